$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for column A (date serials) and column B (values), rows 2-13
$dataA = @(17897, 17928, 17956, 17987, 18017, 18048, 18078, 18109, 18140, 18170, 18201, 18231)
$dataB = @(80.7, 79.1, 81, 83.3, 91, 85.8, 93.4, 96.2, 91.8, 97, 88, 76.9)

for ($i = 0; $i -lt $dataA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dataA[$i]
    $ws.Cells.Item($row, 2).Value = $dataB[$i]
}

# Remove column C entirely (header "Month_num" + all its data)
$ws.Columns.Item(3).Delete()
